$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet (tab name "Checklist" -> "Session")
$ws.Name = "Session"

# Delete row 2 (232002 / Selection / admin@admin.com), shifting rows 3-5 up
$ws.Rows.Item(2).Delete() | Out-Null

# Update the "Type" and "User" columns for the remaining data rows (now rows 2-4)
for ($r = 2; $r -le 4; $r++) {
    $ws.Cells.Item($r, 5).Value = "Scan"
    $ws.Cells.Item($r, 6).Value = "5edfa2692bdacc5e6ee805c626c50cb44cebb065f092d9a1067d89f74dacd326"
}
